# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
#
# Column D ("Price") holds dot-grouped/decimal price strings stored as TEXT in
# the source workbook (e.g. "42.936.24", "236.57"). Excel's COM `Range.Value`
# setter auto-detects numeric-looking text and silently coerces it to a
# Number, which both changes the cell's type and drops formatting such as
# trailing zeros ("236.40" -> 236.4). To keep these cells as text (matching
# the original file) we force the cell's number format to Text ("@") BEFORE
# writing the new value. Column E ("Volume(1h)") values are already padded
# with spaces (e.g. "  +0.75%  "), so Excel stores them as text naturally and
# no extra handling is required there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.959.45"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.297.23"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.23"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.17"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.506"
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("E10").Value = "  +0.84%  "

$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.84"
$ws.Range("E12").Value = "  -2.82%  "

$ws.Range("E13").Value = "  +2.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.04"
$ws.Range("E14").Value = "  +10.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.76"
$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.658.73"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.309.28"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("E18").Value = "  +2.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.893.36"
$ws.Range("E19").Value = "  +0.81%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.61"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.43"
$ws.Range("E23").Value = "  +0.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.40"
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("E25").Value = "  +4.78%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  -1.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.32"
$ws.Range("E28").Value = "  -0.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.64"
$ws.Range("E29").Value = "  +0.32%  "

$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.73"
$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.11"
$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("E34").Value = "  -0.37%  "

$ws.Range("E35").Value = "  +5.64%  "

$ws.Range("E36").Value = "  +1.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.80"
$ws.Range("E37").Value = "  +3.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0701"
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("E43").Value = "  -4.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.990.87"
$ws.Range("E44").Value = "  +1.40%  "

$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("E46").Value = "  +1.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.45"
$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.84"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.528.36"
$ws.Range("E49").Value = "  +0.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.14"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("E51").Value = "  -1.82%  "
